$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.696.48"
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = "'3.788.95"
$ws.Range("E3").Value = '  +1.44%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'595.40"
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = "'3.788.60"
$ws.Range("E7").Value = '  +1.46%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = '  -0.21%  '
$ws.Range("E11").Value = '  -1.92%  '
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").Value = "'36.09"
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("D15").Value = "'4.421.70"
$ws.Range("E15").Value = '  +1.29%  '
$ws.Range("D16").Value = "'3.790.34"
$ws.Range("E16").Value = '  +1.52%  '
$ws.Range("D17").Value = "'18.50"
$ws.Range("E17").Value = '  +3.57%  '
$ws.Range("D18").Value = "'67.693.42"
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").Value = "'7.02"
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("D21").Value = "'10.02"
$ws.Range("E21").Value = '  -6.00%  '
$ws.Range("D22").Value = "'459.23"
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").Value = "'0.697"
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("E24").Value = '  +5.62%  '
$ws.Range("D25").Value = "'83.27"
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("E26").Value = '  +1.38%  '
$ws.Range("D27").Value = "'2.11"
$ws.Range("E27").Value = '  -2.80%  '
$ws.Range("D28").Value = "'10.03"
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  +3.98%  '
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("D33").Value = "'29.66"
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = "'9.08"
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = "'3.727.30"
$ws.Range("E36").Value = '  +1.02%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = "'0.100"
$ws.Range("E37").Value = '  -0.43%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").Value = "'3.39"
$ws.Range("E38").Value = '  -1.36%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = "'0.137"
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").Value = "'0.993"
$ws.Range("E40").Value = '  +0.32%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = "'5.77"
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("B44").Value = 'Arweave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D44").Value = "'44.58"
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = "'48.03"
$ws.Range("E45").Value = '  +2.93%  '
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = "'0.298"
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").Value = "'149.46"
$ws.Range("E47").Value = '  +3.68%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = "'8.29"
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").Value = "'394.83"
$ws.Range("E49").Value = '  +1.80%  '
$ws.Range("E50").Value = '  -4.30%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = "'26.61"
$ws.Range("E51").Value = '  +5.81%  '
